$wb = $excel.ActiveWorkbook

# --- "Per Member Data" sheet: refine averages (E/F/G columns) ---
$pm = $wb.Worksheets.Item("Per Member Data")

$pm.Range("E2").Value = 70.55495978554524
$pm.Range("F2").Value = 34.54700854721421
$pm.Range("G2").Value = 22.764957264931564

$pm.Range("E3").Value = 81.73780487815547
$pm.Range("F3").Value = 46.912408759046556
$pm.Range("G3").Value = 33.394160584141595

$pm.Range("E4").Value = 68.92413793119296
$pm.Range("F4").Value = 26.499999999999996

$pm.Range("E5").Value = 60.54878048791714
$pm.Range("F5").Value = 37.79779411777262
$pm.Range("G5").Value = 22.584558823483277

$pm.Range("E6").Value = 93.25943396213344
$pm.Range("F6").Value = 48.079999999796904
$pm.Range("G6").Value = 41.15999999995648

$pm.Range("E7").Value = 71.02112676060621
$pm.Range("F7").Value = 26.36046511627264
$pm.Range("G7").Value = 28.720930232545275

$pm.Range("E8").Value = 61.00000000000001
$pm.Range("G8").Value = 20.16666666666667

$pm.Range("E9").Value = 89.2399999997534
$pm.Range("F9").Value = 48.079999999796904
$pm.Range("G9").Value = 41.15999999995648

$pm.Range("E10").Value = 73.95032397431406
$pm.Range("F10").Value = 32.95327102811133
$pm.Range("G10").Value = 26.67757009352966

# --- "Drive Team Data" sheet: zero out stale A/B/C columns ---
$dt = $wb.Worksheets.Item("Drive Team Data")

$dt.Range("A2").Value = 0.0
$dt.Range("B2").Value = 0.0
$dt.Range("C2").Value = 0.0

$dt.Range("A3").Value = 0.0
$dt.Range("B3").Value = 0.0
$dt.Range("C3").Value = 0.0

$dt.Range("A4").Value = 0.0
$dt.Range("B4").Value = 0.0
$dt.Range("C4").Value = 0.0
